$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Förändrad" (Changed) date column C, rows 2-238, was bumped by one day
# (serial date 45179 -> 45180) for every record in the sheet.
$lastRow = 238
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45180
